$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.683470726013184
$ws.Range("B1").Value = 1.518706083297729
$ws.Range("C1").Value = 5.441577434539795
$ws.Range("D1").Value = 1.493814468383789
$ws.Range("E1").Value = 0.9092687368392944
